{"js": "// The document contains a single-column table where several rows hold\n// raw benchmark numbers. This edit:\n//  1. Changes row 1's value from \"100\" to \"0M\" and inserts 12 new rows\n//     after it (the GC-log \"0M\"/\"0M\"/\"20\"/... series that used to be\n//     crammed into a later tab-separated row).\n//  2. Changes row 3 (\"70\") to \"0.00000\" and inserts 10 new rows after it\n//     (nine more \"0.00000\" readings plus a trailing \"0.0\").\n//  3. Splits the old tab-separated summary row (\"20<TAB>0.00002<TAB>...\n//     <TAB>100.0\") into individual rows: the row itself becomes empty,\n//     the following (previously empty) row becomes \"100\".\n//  4. Splits the old tab-separated zero row (\"0<TAB>0<TAB>...\") down to\n//     a single \"0\", and appends a new row with \"70\" right after it.\n\nconst table = context.document.body.tables.getFirst();\ntable.rows.load(\"items\");\nawait context.sync();\n\n// Row objects are resolved by their (captured) index, so later inserts\n// shift the position any *earlier*-captured row reference resolves to.\n// Doing the edits from the bottom of the table upward keeps every index\n// used below valid at the time it is used.\nconst rows = table.rows.items;\nconst row1 = rows[0];\nconst row3 = rows[2];\nconst row24 = rows[23];\nconst row25 = rows[24];\nconst row26 = rows[25];\n\nrow1.cells.load(\"items\");\nrow3.cells.load(\"items\");\nrow24.cells.load(\"items\");\nrow25.cells.load(\"items\");\nrow26.cells.load(\"items\");\nawait context.sync();\n\n// --- Step 4 (bottom-most first): collapse the tab-separated zero row to\n// a single \"0\", then append a new row with \"70\" right after it. ---\nrow26.cells.items[0].value = \"0\";\nrow26.insertRows(\"After\", 1, [[\"70\"]]);\nawait context.sync();\n\n// --- Step 3: split the tab-separated summary row + the empty row after it ---\nrow25.cells.items[0].value = \"100\";\nrow24.cells.items[0].value = \"\";\nawait context.sync();\n\n// --- Step 2: row 3 (\"70\" -> \"0.00000\") + insert 10 rows after it ---\nrow3.cells.items[0].value = \"0.00000\";\nrow3.insertRows(\"After\", 10, [\n  [\"0.00000\"],\n  [\"0.00000\"],\n  [\"0.00000\"],\n  [\"0.00000\"],\n  [\"0.00000\"],\n  [\"0.00000\"],\n  [\"0.00000\"],\n  [\"0.00000\"],\n  [\"0.00000\"],\n  [\"0.0\"],\n]);\nawait context.sync();\n\n// --- Step 1: row 1 (\"100\" -> \"0M\") + insert 12 rows after it ---\nrow1.cells.items[0].value = \"0M\";\nrow1.insertRows(\"After\", 12, [\n  [\"0M\"],\n  [\"0M\"],\n  [\"20\"],\n  [\"0.00002\"],\n  [\"0.00004\"],\n  [\"0.00003\"],\n  [\"0.00001\"],\n  [\"0.00003\"],\n  [\"0.00003\"],\n  [\"0.00004\"],\n  [\"0.00063\"],\n  [\"100.0\"],\n]);\nawait context.sync();\n", "ps1": "# The document contains a single-column table where several rows hold\n# raw benchmark numbers. This edit:\n#  1. Changes row 1's value from \"100\" to \"0M\" and inserts 12 new rows\n#     after it (the GC-log \"0M\"/\"0M\"/\"20\"/... series that used to be\n#     crammed into a later tab-separated row).\n#  2. Changes row 3 (\"70\") to \"0.00000\" and inserts 10 new rows after it\n#     (nine more \"0.00000\" readings plus a trailing \"0.0\").\n#  3. Splits the old tab-separated summary row (\"20<TAB>0.00002<TAB>...\n#     <TAB>100.0\") into individual rows: the row itself becomes empty,\n#     the following (previously empty) row becomes \"100\".\n#  4. Splits the old tab-separated zero row (\"0<TAB>0<TAB>...\") down to\n#     a single \"0\", and appends a new row with \"70\" right after it.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Row objects resolve by their (live) index, so earlier-captured\n# references go stale once rows are inserted above them. Doing the\n# edits from the bottom of the table upward keeps every index used\n# below valid at the time it is used.\n\n# --- Step 4 (bottom-most first): collapse the tab-separated zero row to\n# a single \"0\", then append a new row with \"70\" right after it. ---\n$row26 = $t.Rows.Item(26)\n$row26.Cells.Item(1).Range.Text = \"0\"\n$newRow = $t.Rows.Add()\n$newRow.Cells.Item(1).Range.Text = \"70\"\n\n# --- Step 3: split the tab-separated summary row + the empty row after it ---\n$row25 = $t.Rows.Item(25)\n$row25.Cells.Item(1).Range.Text = \"100\"\n$row24 = $t.Rows.Item(24)\n$row24.Cells.Item(1).Range.Text = \"\"\n\n# --- Step 2: row 3 (\"70\" -> \"0.00000\") + insert 10 rows after it ---\n$row3 = $t.Rows.Item(3)\n$row3.Cells.Item(1).Range.Text = \"0.00000\"\n$row2Vals = @(\"0.00000\",\"0.00000\",\"0.00000\",\"0.00000\",\"0.00000\",\"0.00000\",\"0.00000\",\"0.00000\",\"0.00000\",\"0.0\")\n$beforeRow = $t.Rows.Item(4)\nfor ($i = $row2Vals.Length - 1; $i -ge 0; $i--) {\n  $newRow = $t.Rows.Add($beforeRow)\n  $newRow.Cells.Item(1).Range.Text = $row2Vals[$i]\n}\n\n# --- Step 1: row 1 (\"100\" -> \"0M\") + insert 12 rows after it ---\n$row1 = $t.Rows.Item(1)\n$row1.Cells.Item(1).Range.Text = \"0M\"\n$row1Vals = @(\"0M\",\"0M\",\"20\",\"0.00002\",\"0.00004\",\"0.00003\",\"0.00001\",\"0.00003\",\"0.00003\",\"0.00004\",\"0.00063\",\"100.0\")\n$beforeRow = $t.Rows.Item(2)\nfor ($i = $row1Vals.Length - 1; $i -ge 0; $i--) {\n  $newRow = $t.Rows.Add($beforeRow)\n  $newRow.Cells.Item(1).Range.Text = $row1Vals[$i]\n}\n"}
